$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Rewrite the author-list sentence:
#    "Nair, G.S., Bhat, C.R., Pendyala, R.M., Loo, B.P.Y. and Lam, W.H.K., "
#    -> "Nair, G.S., C.R. Bhat, R.M. Pendyala, B.P.Y. Loo and W.H.K. Lam, "
#    split across many runs (matching the target XML), with a "_GoBack"
#    bookmark inserted right before the final "Lam" run.
# ---------------------------------------------------------------------------

$oldSentence = "Nair, G.S., Bhat, C.R., Pendyala, R.M., Loo, B.P.Y. and Lam, W.H.K., "
$newSentence = "Nair, G.S., C.R. Bhat, R.M. Pendyala, B.P.Y. Loo and W.H.K. Lam, "

$findRng = $d.Content
$found = $findRng.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the author sentence to rewrite"
}
$start = $findRng.Start

# Replace the whole sentence with the new text first (single run for now).
$wholeRng = $d.Range($start, $findRng.End)
$wholeRng.Text = $newSentence

# Pieces the new sentence should be broken into (as separate <w:r> runs).
# The 13th boundary (index 12, right after "W.H.K. ") is where the _GoBack
# bookmark belongs, before the final "Lam" / ", " runs.
$pieces = @(
    "Nair, G.S., ",
    "C.R.",
    " Bhat",
    ", ",
    "R.M.",
    " ",
    "Pendy",
    "ala, ",
    "B.P.Y.",
    " Loo ",
    "and ",
    "W.H.K.",
    " ",
    "Lam",
    ", "
)

# Compute absolute offsets of every boundary between pieces (not counting
# the very first / very last edges, since those already coincide with the
# surrounding text and need no wedge).
$boundaries = @()
$pos = $start
for ($i = 0; $i -lt ($pieces.Count - 1); $i++) {
    $pos = $pos + $pieces[$i].Length
    $boundaries += $pos
}

# Index (within $boundaries) of the boundary where _GoBack belongs: right
# after piece 12 (" ") and before piece 13 ("Lam") -> that's boundary index 12.
$goBackBoundaryIndex = 12

# Insert a temporary "wedge" bookmark at every internal boundary, from the
# last boundary back to the first, so runs split apart into separate <w:r>
# elements instead of being coalesced back together.
for ($i = $boundaries.Count - 1; $i -ge 0; $i--) {
    $p = $boundaries[$i]
    $wedgeRange = $d.Range($p, $p)
    $d.Bookmarks.Add("ZZWedge$i", $wedgeRange)
}

# Drop every wedge bookmark again (the run split persists even after the
# bookmark that caused it is gone) except for the one where _GoBack should
# live; turn that one into the real _GoBack bookmark.
for ($i = 0; $i -lt $boundaries.Count; $i++) {
    if ($i -ne $goBackBoundaryIndex) {
        $d.Bookmarks("ZZWedge$i").Delete()
    }
}

$goBackPos = $boundaries[$goBackBoundaryIndex]
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
if ($d.Bookmarks.Exists("ZZWedge$goBackBoundaryIndex")) {
    $d.Bookmarks("ZZWedge$goBackBoundaryIndex").Delete()
}

# ---------------------------------------------------------------------------
# 2) The old "_GoBack" bookmark (after "(dependent variables) ") is removed.
#    Adding the bookmark above already relocated it away from that spot
#    (a document can only have one bookmark with a given name), but make
#    sure no stray copy is left behind just in case.
# ---------------------------------------------------------------------------
# (No further action required: Bookmarks.Add with an existing name moves it.)
